$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 88-92: new ring oscillator / timing control registers ---
# (order chosen to reproduce the original shared-string insertion order)
$ws.Range("D89").Value = "LSB"
$ws.Range("C89").Value = "ENABLE FIRMWARE FEEDBACK LOOP for ROvcp"

$ws.Range("D90").Value = "LSB"

# --- Rows 11/12: register 0x04/0x05 area ---
$ws.Range("D11").Value = "read_only"
$ws.Range("C12").Value = "ring oscillator counter divider output"
$ws.Range("D12").Value = "read_only"

$ws.Range("C90").Value = "trig sign - in psec4a serial block"

$ws.Range("D91").Value = "LSB"
$ws.Range("C91").Value = "dll speed select - in psec4a serial block"

$ws.Range("D92").Clear()
$ws.Range("D92").Value = "LSB"
$ws.Range("C92").Value = "use reset in xfer process - in psec4a serial block"

$ws.Range("C88").Value = "COUNT TARGET for ring osc when reg 0x52 enabled"
$ws.Range("D88").Clear()

# --- Update the active selection to match the edited workbook's last position ---
$ws.Range("C90").Select()
